$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (AUTO, 1090+)
$ws.Range("E2").Value = 186
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 60058
$ws.Range("H2").Value = 0.17

# Row 3 (AUTO, 366 TO 730)
$ws.Range("E3").Value = 114
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = 850757.8
$ws.Range("H3").Value = 7.68

# Row 4 (AUTO, 731 TO 1090)
$ws.Range("E4").Value = 43
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 408501
$ws.Range("H4").Value = 7.65
